$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 255.83333
$ws.Range("I6").Value = 91.8
$ws.Range("J6").Value = 299
$ws.Range("K6").Value = 275.4
$ws.Range("L6").Value = 897
$ws.Range("M6").Value = -163.4
$ws.Range("N6").Value = -1121

$ws.Range("H17").Value = 885.2133
$ws.Range("J17").Value = 891.1081
$ws.Range("L17").Value = 2673.3243
$ws.Range("N17").Value = -3009.3243

$ws.Range("H40").Value = 10004930
$ws.Range("I40").Value = 25002350
$ws.Range("J40").Value = 6650
$ws.Range("K40").Value = 25002350
$ws.Range("L40").Value = 6650
$ws.Range("M40").Value = -25002175
$ws.Range("N40").Value = -7000

$ws.Range("H62").Value = 201622
$ws.Range("I62").Value = 201622
$ws.Range("K62").Value = 201622
$ws.Range("M62").Value = -200998

$ws.Range("H65").Value = 201622
$ws.Range("I65").Value = 201622
$ws.Range("K65").Value = 1008110
$ws.Range("M65").Value = -1004990

$ws.Range("H100").Value = 9866.959999999999
$ws.Range("I100").Value = 4707
$ws.Range("K100").Value = 4707
$ws.Range("M100").Value = -4166

$ws.Range("H112").Value = 669398.9
$ws.Range("J112").Value = 717070.9399999999
$ws.Range("L112").Value = 2151212.82
$ws.Range("N112").Value = -2153428.82

$ws.Range("H116").Value = 2790.375
$ws.Range("I116").Value = 2584.8
$ws.Range("K116").Value = 2584.8
$ws.Range("M116").Value = 857.1999999999998

$ws.Range("H132").Value = 4841.415
$ws.Range("I132").Value = 4902.5713
$ws.Range("K132").Value = 14707.7139
$ws.Range("M132").Value = -12177.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1819.1111
$ws.Range("I32").Value = 1876.5072
$ws.Range("K32").Value = 1876.5072
$ws.Range("M32").Value = -1589.5072

$ws.Range("H61").Value = 5679.8687
$ws.Range("I61").Value = 3640.739
$ws.Range("K61").Value = 3640.739
$ws.Range("M61").Value = -3428.739

$ws.Range("H74").Value = 1204
$ws.Range("I74").Value = 1133.579
$ws.Range("K74").Value = 1133.579
$ws.Range("M74").Value = -259.579

$ws.Range("H76").Value = 60963.832
$ws.Range("I76").Value = 19995
$ws.Range("K76").Value = 19995
$ws.Range("M76").Value = -19657

$ws.Range("H77").Value = 1204
$ws.Range("I77").Value = 1133.579
$ws.Range("K77").Value = 5667.895
$ws.Range("M77").Value = -1299.895

$ws.Range("H79").Value = 60963.832
$ws.Range("I79").Value = 19995
$ws.Range("K79").Value = 19995
$ws.Range("M79").Value = -18825

$ws.Range("H97").Value = 1362.027
$ws.Range("I97").Value = 852
$ws.Range("K97").Value = 852
$ws.Range("M97").Value = -356

$ws.Range("H105").Value = 98869.5
$ws.Range("J105").Value = 98869.5
$ws.Range("L105").Value = 98869.5
$ws.Range("N105").Value = -105857.5

$ws.Range("H122").Value = 11639.857
$ws.Range("I122").Value = 11530.733
$ws.Range("J122").Value = 11912.667
$ws.Range("K122").Value = 34592.199
$ws.Range("L122").Value = 35738.001
$ws.Range("M122").Value = -32142.199
$ws.Range("N122").Value = -40638.001

$ws.Range("H132").Value = 3637.7964
$ws.Range("I132").Value = 1604.0488
$ws.Range("K132").Value = 4812.1464
$ws.Range("M132").Value = -2282.1464

$ws.Range("H136").Value = 5679.8687
$ws.Range("I136").Value = 3640.739
$ws.Range("K136").Value = 10922.217
$ws.Range("M136").Value = -8372.217000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 842.8570999999999
$ws.Range("I22").Value = 733.3333
$ws.Range("K22").Value = 733.3333
$ws.Range("M22").Value = -560.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4322
$ws.Range("I16").Value = 3807.6538
$ws.Range("J16").Value = 6996.6
$ws.Range("K16").Value = 3807.6538
$ws.Range("L16").Value = 6996.6
$ws.Range("M16").Value = -3520.6538
$ws.Range("N16").Value = -7570.6

$ws.Range("H31").Value = 2920.9092
$ws.Range("I31").Value = 1348.7333
$ws.Range("K31").Value = 1348.7333
$ws.Range("M31").Value = -1053.7333

$ws.Range("H34").Value = 2920.9092
$ws.Range("I34").Value = 1348.7333
$ws.Range("K34").Value = 1348.7333
$ws.Range("M34").Value = -1146.7333

$ws.Range("H105").Value = 4915.12
$ws.Range("I105").Value = 3243.7778
$ws.Range("K105").Value = 3243.7778
$ws.Range("M105").Value = -1496.7778

$ws.Range("H113").Value = 4322
$ws.Range("I113").Value = 3807.6538
$ws.Range("J113").Value = 6996.6
$ws.Range("K113").Value = 3807.6538
$ws.Range("L113").Value = 6996.6
$ws.Range("M113").Value = -1637.6538
$ws.Range("N113").Value = -11336.6

$ws.Range("H132").Value = 1325.3077
$ws.Range("I132").Value = 1289.9166
$ws.Range("J132").Value = 1750
$ws.Range("K132").Value = 3869.7498
$ws.Range("L132").Value = 5250
$ws.Range("M132").Value = -1339.7498
$ws.Range("N132").Value = -10310

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 161.5625
$ws.Range("I2").Value = 188.38461
$ws.Range("J2").Value = 45.333332
$ws.Range("K2").Value = 188.38461
$ws.Range("L2").Value = 45.333332
$ws.Range("M2").Value = -75.38461000000001
$ws.Range("N2").Value = -271.333332

$ws.Range("H59").Value = 10000
$ws.Range("J59").Value = 10000
$ws.Range("L59").Value = 10000
$ws.Range("N59").Value = -11166

$ws.Range("H113").Value = 5269.7144
$ws.Range("I113").Value = 6192.4
$ws.Range("K113").Value = 6192.4
$ws.Range("M113").Value = -4022.4

$ws.Range("H122").Value = 3451.625
$ws.Range("I122").Value = 4118.8335
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 12356.5005
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -9906.500499999998
$ws.Range("N122").Value = -9250

$ws.Range("H132").Value = 10806.586
$ws.Range("I132").Value = 10755.259
$ws.Range("K132").Value = 32265.777
$ws.Range("M132").Value = -29735.777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 28586.477
$ws.Range("I122").Value = 29126.188
$ws.Range("J122").Value = 26859.4
$ws.Range("K122").Value = 87378.564
$ws.Range("L122").Value = 80578.20000000001
$ws.Range("M122").Value = -84928.564
$ws.Range("N122").Value = -85478.20000000001

$ws.Range("H136").Value = 4749.2974
$ws.Range("I136").Value = 4756.6177
$ws.Range("K136").Value = 14269.8531
$ws.Range("M136").Value = -11719.8531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 4013
$ws.Range("I34").Value = 3026
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 3026
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -2823
$ws.Range("N34").Value = -5406

$ws.Range("H54").Value = 27639.545
$ws.Range("I54").Value = 13035
$ws.Range("K54").Value = 13035
$ws.Range("M54").Value = -12515

$ws.Range("H81").Value = 4481.4707
$ws.Range("I81").Value = 2198.6667
$ws.Range("J81").Value = 7049.625
$ws.Range("K81").Value = 4397.3334
$ws.Range("L81").Value = 14099.25
$ws.Range("M81").Value = -3336.3334
$ws.Range("N81").Value = -16221.25

$ws.Range("H84").Value = 4481.4707
$ws.Range("I84").Value = 2198.6667
$ws.Range("J84").Value = 7049.625
$ws.Range("K84").Value = 21986.667
$ws.Range("L84").Value = 70496.25
$ws.Range("M84").Value = -16682.667
$ws.Range("N84").Value = -81104.25

$ws.Range("H107").Value = 482.87878
$ws.Range("I107").Value = 269.6
$ws.Range("J107").Value = 811
$ws.Range("K107").Value = 808.8000000000001
$ws.Range("L107").Value = 2433
$ws.Range("M107").Value = 1111.2
$ws.Range("N107").Value = -6273

$ws.Range("H122").Value = 2601.4167
$ws.Range("I122").Value = 2771.9473
$ws.Range("K122").Value = 8315.841899999999
$ws.Range("M122").Value = -5865.841899999999

$ws.Range("H126").Value = 10730.8
$ws.Range("I126").Value = 9154.058999999999
$ws.Range("K126").Value = 27462.177
$ws.Range("M126").Value = -24992.177

$ws.Range("H132").Value = 1624.52
$ws.Range("I132").Value = 1470.75
$ws.Range("J132").Value = 2239.6
$ws.Range("K132").Value = 4412.25
$ws.Range("L132").Value = 6718.799999999999
$ws.Range("M132").Value = -1882.25
$ws.Range("N132").Value = -11778.8
